$d = $word.ActiveDocument

# --- 1. "Backtracking:" run gets cyan highlight ---
$r1 = $d.Content
$null = $r1.Find.Execute("Backtracking:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Font.HighlightColorIndex = 3

# --- 2. Split " N-Queens problem, Sudoku solver, Rat in a maze" so that
#        " N-Queens problem," becomes its own (cyan-highlighted) run and
#        " Sudoku solver, Rat in a maze" remains a separate, unhighlighted run ---
$r2 = $d.Content
$null = $r2.Find.Execute(" N-Queens problem,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Font.HighlightColorIndex = 3

# --- 3. "Recursion vs. Iteration" paragraph: highlight the run AND the
#        paragraph mark (pPr/rPr) yellow, by applying to the whole paragraph
#        range (which includes the trailing paragraph mark) ---
$r3 = $d.Content
$null = $r3.Find.Execute("Recursion vs. Iteration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p3 = $r3.Paragraphs(1)
$pr3 = $p3.Range
$pr3.Font.HighlightColorIndex = 7
